$wb = $excel.ActiveWorkbook

# 1) Rename "test1" -> "weizmann"
$weizmann = $wb.Worksheets.Item("test1")
$weizmann.Name = "weizmann"

# 2) Insert a brand-new sheet right after "weizmann" (i.e. before "test3"),
#    named "weizmann_testing"
$weizmannTesting = $wb.Worksheets.Add($null, $weizmann)
$weizmannTesting.Name = "weizmann_testing"

# 3) weizmann_testing gets the same header row as weizmann plus two data
#    rows for the new "lyova" / "moshe" subjects
$weizmannTesting.Range("A1").Value = "RUN"
$weizmannTesting.Range("B1").Value = "WALK"
$weizmannTesting.Range("C1").Value = "WAVE2"
$weizmannTesting.Range("D1").Value = "JUMP"
$weizmannTesting.Range("E1").Value = "PJUMP"
$weizmannTesting.Range("F1").Value = "JACK"
$weizmannTesting.Range("G1").Value = "SIDE"
$weizmannTesting.Range("H1").Value = "SKIP"
$weizmannTesting.Range("I1").Value = "WAVE1"

$weizmannTesting.Range("A2").Value = "lyova_run.avi"
$weizmannTesting.Range("B2").Value = "lyova_walk.avi"
$weizmannTesting.Range("C2").Value = "lyova_wave2.avi"
$weizmannTesting.Range("D2").Value = "lyova_jump.avi"
$weizmannTesting.Range("E2").Value = "lyova_pjump.avi"
$weizmannTesting.Range("F2").Value = "lyova_jack.avi"
$weizmannTesting.Range("G2").Value = "lyova_side.avi"
$weizmannTesting.Range("H2").Value = "lyova_skip.avi"
$weizmannTesting.Range("I2").Value = "lyova_wave1.avi"

$weizmannTesting.Range("A3").Value = "moshe_run.avi"
$weizmannTesting.Range("B3").Value = "moshe_walk.avi"
$weizmannTesting.Range("C3").Value = "moshe_wave2.avi"
$weizmannTesting.Range("D3").Value = "moshe_jump.avi"
$weizmannTesting.Range("E3").Value = "moshe_pjump.avi"
$weizmannTesting.Range("F3").Value = "moshe_jack.avi"
$weizmannTesting.Range("G3").Value = "moshe_side.avi"
$weizmannTesting.Range("H3").Value = "moshe_skip.avi"
$weizmannTesting.Range("I3").Value = "moshe_wave1.avi"

# 4) weizmann row 8 used to be a near-empty "lena2" stub row; it is now
#    filled in with real "shahar" data
$weizmann.Range("A8").Value = "shahar_run.avi"
$weizmann.Range("B8").Value = "shahar_walk.avi"
$weizmann.Range("C8").Value = "shahar_wave2.avi"
$weizmann.Range("D8").Value = "shahar_jump.avi"
$weizmann.Range("E8").Value = "shahar_pjump.avi"
$weizmann.Range("F8").Value = "shahar_jack.avi"
$weizmann.Range("G8").Value = "shahar_side.avi"
$weizmann.Range("H8").Value = "shahar_skip.avi"
$weizmann.Range("I8").Value = "shahar_wave1.avi"

# 5) Update selections to match the new state
[void]$weizmann.Range("F8").Select()
[void]$weizmannTesting.Range("I6").Select()

# "pagi" had the keyboard focus before; the new "weizmann_testing" sheet
# is the active one now
[void]$weizmannTesting.Activate()
